$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (risk 2): rename risk description, add mitigation text ---
$ws.Range("C4").Value = "Desfalque de algum membro do grupo"

# --- Row 8 (risk 6): fill in the previously empty risk row ---
$ws.Range("C8").Value = "Projeto desorganizado"

# --- Row 8 H: action text ---
$ws.Range("H8").Value = "Projeto dividido e organizado nas ferramentas de gestão ,juntamente com daily's"

# --- Row 4 H: action text ---
$ws.Range("H4").Value = "Manter os integrantes engajados com o projeto"

# --- Row 5 H: action text ---
$ws.Range("H5").Value = "Aumentar o fluxo de comunicação entre o grupo e flexibilização nos horários das daily's para evitar ausências"

# --- Row 6 H: action text ---
$ws.Range("H6").Value = "Ter a disponibilidade de frequentar diversos pontos de conexão a rede"

# --- Row 9 (risk 7): fill in the previously empty risk row ---
$ws.Range("C9").Value = "Sprints desbalanceadas"

# --- Row 9 H: action text ---
$ws.Range("H9").Value = "Fazendo as estimatívas dos tamanhos dos requisitos e criando o gráfico de Burndown "

# --- Fill remaining values for rows 5, 8, 9 ---
$ws.Range("G5").Value = "Evitar"
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = "Evitar"
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = "Evitar"

# --- New rows 10 and 11 (risks 8 and 9) ---
$ws.Range("B10").Value = 8
$ws.Range("B11").Value = 9
